$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(471).Insert()

$ws.Cells.Item(471, 1).Value = 7
$ws.Cells.Item(471, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(471, 3).Value = "Ñuble"
$ws.Cells.Item(471, 4).Value = 45239
$ws.Cells.Item(471, 5).Value = 16
$ws.Cells.Item(471, 6).Value = "Fruta"
$ws.Cells.Item(471, 7).Value = 100102
$ws.Cells.Item(471, 8).Value = "Cítricos"
$ws.Cells.Item(471, 9).Value = 100102004
$ws.Cells.Item(471, 10).Value = "Mandarina"
$ws.Cells.Item(471, 11).Value = "Clementina"
$ws.Cells.Item(471, 12).Value = "Primera"
$ws.Cells.Item(471, 13).Value = 100
$ws.Cells.Item(471, 14).Value = 10000
$ws.Cells.Item(471, 15).Value = 10000
$ws.Cells.Item(471, 16).Value = 10000
$ws.Cells.Item(471, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(471, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(471, 19).Value = 556
$ws.Cells.Item(471, 20).Value = 18
